$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '98.447.42'
$ws.Range('E2').Value = '  +1.02%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.491.74'
$ws.Range('E3').Value = '  +4.08%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '253.13'
$ws.Range('E5').Value = '  +1.39%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '666.37'
$ws.Range('E6').Value = '  +1.73%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.48'
$ws.Range('E7').Value = '  +6.33%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.428'
$ws.Range('E8').Value = '  +1.95%  '
$ws.Range('E9').Value = '  +3.45%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.00'
$ws.Range('E10').Value = '  +0.07%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '3.490.37'
$ws.Range('E11').Value = '  +4.20%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '45.64'
$ws.Range('E12').Value = '  +12.58%  '
$ws.Range('E13').Value = '  +1.02%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.27'
$ws.Range('E14').Value = '  +3.09%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '98.267.69'
$ws.Range('E15').Value = '  +0.97%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000261'
$ws.Range('E16').Value = '  +2.34%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.151.71'
$ws.Range('E17').Value = '  +4.43%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '8.98'
$ws.Range('E18').Value = '  +3.68%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.493.97'
$ws.Range('E19').Value = '  +4.31%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.77'
$ws.Range('E20').Value = '  +11.25%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.75'
$ws.Range('E21').Value = '  +9.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.528'
$ws.Range('E22').Value = '  -4.72%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '524.54'
$ws.Range('E23').Value = '  +4.06%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.43'
$ws.Range('E24').Value = '  +2.27%  '
$ws.Range('E25').Value = '  +1.34%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.80'
$ws.Range('E26').Value = '  +7.45%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '98.52'
$ws.Range('E27').Value = '  +1.83%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.77'
$ws.Range('E28').Value = '  +5.14%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.682.26'
$ws.Range('E29').Value = '  +4.11%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '12.44'
$ws.Range('E30').Value = '  +12.21%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.87'
$ws.Range('E31').Value = '  +14.00%  '
$ws.Range('E32').Value = '  -1.34%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.995'
$ws.Range('E33').Value = '  -0.01%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.190'
$ws.Range('E34').Value = '  -0.33%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.598'
$ws.Range('E35').Value = '  +8.33%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '31.25'
$ws.Range('E36').Value = '  +9.10%  '
$ws.Range('E37').Value = '  +0.55%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.52'
$ws.Range('E38').Value = '  +3.59%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '7.97'
$ws.Range('E39').Value = '  +3.02%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.156'
$ws.Range('E40').Value = '  +4.84%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '526.95'
$ws.Range('E41').Value = '  +1.60%  '
$ws.Range('B42').Value = 'USDe'
$ws.Range('C42').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.08%  '
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.919'
$ws.Range('E43').Value = '  +8.62%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.77'
$ws.Range('E44').Value = '  +7.13%  '
$ws.Range('B45').Value = 'WhiteBITCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '24.42'
$ws.Range('E45').Value = '  -0.85%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0435'
$ws.Range('E46').Value = '  +3.40%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.79'
$ws.Range('E47').Value = '  +3.69%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.73'
$ws.Range('E48').Value = '  -1.25%  '
$ws.Range('E49').Value = '  -1.04%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.23'
$ws.Range('E50').Value = '  +11.62%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '55.58'
$ws.Range('E51').Value = '  +4.15%  '
